$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DDT - the departure_date column should hold plain text dates (so the
# values can be read verbatim the way they'd come out of a CSV file)
# instead of native Excel date serials.
$ws.Range("C2:C4").NumberFormat = "@"

$ws.Range("C2").Value = "10/15/2023"
$ws.Range("C3").Value = "10/20/2023"
$ws.Range("C4").Value = "11/03/2023"

# Leave the selection where the author last left it.
$ws.Range("D10").Select()
